$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wrap the "écart" (difference) formulas in ABS() so the deviation is
# always reported as a positive value.
$cols = @("C", "D", "E", "F", "G", "H")
foreach ($col in $cols) {
    $ws.Range($col + "10").Formula = "=ABS(" + $col + "2-" + $col + "4)"
    $ws.Range($col + "11").Formula = "=ABS(" + $col + "3-" + $col + "5)"
}

# Thin out the thick outer border that used to close off the bottom-right
# of the table (right edge of column H, bottom edge of row 11).
$rng = $ws.Range("C10:H11")
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2

# Leave the selection where the author left it.
$ws.Range("K10").Select() | Out-Null
